$wb = $excel.ActiveWorkbook

# 展览 (sheet1) - update F column "想去人数" (want-to-go count) values
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 570
$ws.Range("F3").Value = 5367
$ws.Range("F8").Value = 372
$ws.Range("F9").Value = 1338
$ws.Range("F12").Value = 3060
$ws.Range("F13").Value = 1902
$ws.Range("F15").Value = 56
$ws.Range("F17").Value = 14
$ws.Range("F20").Value = 968
$ws.Range("F23").Value = 3507
$ws.Range("F24").Value = 1102
$ws.Range("F25").Value = 2788
$ws.Range("F27").Value = 1841
$ws.Range("F28").Value = 4018
$ws.Range("F29").Value = 106
$ws.Range("F30").Value = 912
$ws.Range("F31").Value = 458
$ws.Range("F33").Value = 5
$ws.Range("F34").Value = 26
$ws.Range("F35").Value = 990
$ws.Range("F36").Value = 1256
$ws.Range("F38").Value = 1019
$ws.Range("F39").Value = 661
$ws.Range("F40").Value = 510

# 演出 (sheet2) - update F column "想去人数" (want-to-go count) values
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 903

# 全部类型 (sheet4) - update F column "想去人数" (want-to-go count) values
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 570
$ws.Range("F3").Value = 570
$ws.Range("F4").Value = 5367
$ws.Range("F9").Value = 372
$ws.Range("F10").Value = 1338
$ws.Range("F11").Value = 3060
$ws.Range("F13").Value = 1902
$ws.Range("F15").Value = 56
$ws.Range("F18").Value = 903
$ws.Range("F21").Value = 132
$ws.Range("F22").Value = 968
$ws.Range("F24").Value = 3507
$ws.Range("F27").Value = 1102
$ws.Range("F28").Value = 2788
$ws.Range("F29").Value = 1844
$ws.Range("F30").Value = 4018
$ws.Range("F32").Value = 106
$ws.Range("F33").Value = 912
$ws.Range("F35").Value = 26
$ws.Range("F36").Value = 990
$ws.Range("F38").Value = 1256
$ws.Range("F40").Value = 1019
$ws.Range("F42").Value = 661

